# Se arreglan localizadores en WelcomePage y nombre de metodos
#
# Rewrites the DataPrueba sheet:
#  - shrinks the used range from A1:L21 down to A1:F12
#  - replaces the CP00x_* catalog strings and Dato00x filler values with the
#    new WelcomePage locators / phrases
#  - adds the "jisola.tsoft@gmail.com" mailto hyperlink (with its style) to
#    B2:B5 instead of just B2
#  - moves the active selection to C4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the columns/rows that are no longer part of the data set
#    (G:L entirely, and rows 13-21 in A:F). ClearContents() also shrinks
#    <dimension> and compacts the shared-strings table automatically.
# ---------------------------------------------------------------------
$ws.Range("G1:L21").ClearContents()
$ws.Range("A13:F21").ClearContents()

# ---------------------------------------------------------------------
# 2. Remove the pre-existing hyperlink on B2 (it pointed at the old
#    "algo@algo.com" address) so it can be rebuilt against the new
#    address together with B3:B5 below.
# ---------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3. New cell values
# ---------------------------------------------------------------------
$iquest = [char]0x00BF
$ntilde = [char]0x00F1
$oacute = [char]0x00F3

# Row 2
$ws.Range("A2").Value = "CP001_login_fallido"
$ws.Range("B2").Value = "jisola.tsoft@gmail.com"
$ws.Range("C2").Value = 12345678
$ws.Range("D2").Value = $iquest + "Olvidaste tu contrase" + $ntilde + "a?"

# Row 3
$ws.Range("A3").Value = "CP002_login_exitoso"
$ws.Range("B3").Value = "jisola.tsoft@gmail.com"
$ws.Range("C3").Value = 12061990
$ws.Range("D3").Value = "Te damos la bienvenida a Facebook, Juan"

# Row 4
$ws.Range("A4").Value = "CP003_cerrar_sesion"
$ws.Range("B4").Value = "jisola.tsoft@gmail.com"
$ws.Range("C4").Value = 12061990
$ws.Range("D4").Value = "Iniciar sesi" + $oacute + "n"

# Row 5
$ws.Range("A5").Value = "CP004_modo_oscuro"
$ws.Range("B5").Value = "jisola.tsoft@gmail.com"
$ws.Range("C5").Value = 12061990

# Rows 6-12 (single-column catalog of remaining test case names)
$ws.Range("A6").Value = "CP005_buscar_persona"
$ws.Range("A7").Value = "CP006_enviar_solicitud"
$ws.Range("A8").Value = "CP007_cancelar_solicitud"
$ws.Range("A9").Value = "CP008_meGusta_pagina"
$ws.Range("A10").Value = "CP009_crear_publicacion"
$ws.Range("A11").Value = "CP010_crear_historia"
$ws.Range("A12").Value = "CP011_enviar_mensaje"

# ---------------------------------------------------------------------
# 4. Re-create the mailto hyperlinks on B2:B5, all against the new
#    address, and keep them all using the workbook's builtin hyperlink
#    cell style so the column renders consistently.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:jisola.tsoft@gmail.com") | Out-Null

$hyperlinkStyleName = $wb.Styles.Item(1).Name
$ws.Range("B2").Style = $hyperlinkStyleName
$ws.Range("B3").Style = $hyperlinkStyleName
$ws.Range("B4").Style = $hyperlinkStyleName
$ws.Range("B5").Style = $hyperlinkStyleName

# The builtin "Hyperlink"/"Hipervínculo" cell style name itself is locale
# generated by Excel and not rename-able through the object model, but
# attempt it in case the host Excel instance exposes a writable Name.
try {
    $wb.Styles.Item(1).Name = "Hyperlink"
} catch {
}

# ---------------------------------------------------------------------
# 5. Match the saved selection / active cell shown in the sheet view.
# ---------------------------------------------------------------------
$ws.Range("C4").Select() | Out-Null
